# Weekly update: add a new day's (2021-10-26, serial 44495) price records
# for Feria Lagunitas de Puerto Montt - Pera, inserted at the top of the
# existing data block (original rows 60-157 shift down to 63-160).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows right before the current first data row of the new
# date range (row 60), pushing all existing data (old rows 60-157) down
# to rows 63-160.
$ws.Rows("60:62").Insert()

# --- New row 60: Pera / Forelle / Primera ---
$ws.Range("A60").Value = 4
$ws.Range("B60").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C60").Value = "Los Lagos"
$ws.Range("D60").Value = 44495
$ws.Range("E60").Value = 10
$ws.Range("F60").Value = "Fruta"
$ws.Range("G60").Value = 100104
$ws.Range("H60").Value = "Frutos de pepita"
$ws.Range("I60").Value = 100104005
$ws.Range("J60").Value = "Pera"
$ws.Range("K60").Value = "Forelle"
$ws.Range("L60").Value = "Primera"
$ws.Range("M60").Value = 500
$ws.Range("N60").Value = 13000
$ws.Range("O60").Value = 13500
$ws.Range("P60").Value = 13250
$ws.Range("Q60").Value = "$/caja 15 kilos empedrada"
$ws.Range("R60").Value = "Región de O'Higgins"
$ws.Range("S60").Value = 883
$ws.Range("T60").Value = 15

# --- New row 61: Pera / Packham's Triumph / Primera ---
$ws.Range("A61").Value = 4
$ws.Range("B61").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C61").Value = "Los Lagos"
$ws.Range("D61").Value = 44495
$ws.Range("E61").Value = 10
$ws.Range("F61").Value = "Fruta"
$ws.Range("G61").Value = 100104
$ws.Range("H61").Value = "Frutos de pepita"
$ws.Range("I61").Value = 100104005
$ws.Range("J61").Value = "Pera"
$ws.Range("K61").Value = "Packham's Triumph"
$ws.Range("L61").Value = "Primera"
$ws.Range("M61").Value = 500
$ws.Range("N61").Value = 15000
$ws.Range("O61").Value = 16000
$ws.Range("P61").Value = 15500
$ws.Range("Q61").Value = "$/caja 15 kilos empedrada"
$ws.Range("R61").Value = "Región de O'Higgins"
$ws.Range("S61").Value = 1033
$ws.Range("T61").Value = 15

# --- New row 62: Pera / Packham's Triumph / Segunda ---
$ws.Range("A62").Value = 4
$ws.Range("B62").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C62").Value = "Los Lagos"
$ws.Range("D62").Value = 44495
$ws.Range("E62").Value = 10
$ws.Range("F62").Value = "Fruta"
$ws.Range("G62").Value = 100104
$ws.Range("H62").Value = "Frutos de pepita"
$ws.Range("I62").Value = 100104005
$ws.Range("J62").Value = "Pera"
$ws.Range("K62").Value = "Packham's Triumph"
$ws.Range("L62").Value = "Segunda"
$ws.Range("M62").Value = 200
$ws.Range("N62").Value = 13000
$ws.Range("O62").Value = 13000
$ws.Range("P62").Value = 13000
$ws.Range("Q62").Value = "$/caja 15 kilos empedrada"
$ws.Range("R62").Value = "Región de O'Higgins"
$ws.Range("S62").Value = 867
$ws.Range("T62").Value = 15
